# Fruta / hortaliza, semanal
# Re-shuffle the D (Fecha), J (Volumen), K (Precio mínimo), L (Precio máximo),
# M (Precio promedio ponderado) and P (Precio $/Kg) values across data rows
# 2-27, per the row permutation captured from the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Maps destination row -> source row (the values that should end up in the
# destination row are the *original* values found in the source row).
$rowMap = @{
    2  = 3
    3  = 14
    4  = 16
    5  = 8
    6  = 13
    7  = 12
    8  = 7
    9  = 23
    10 = 27
    11 = 17
    12 = 26
    13 = 4
    14 = 25
    15 = 24
    16 = 9
    17 = 22
    18 = 10
    19 = 5
    20 = 2
    21 = 6
    22 = 18
    23 = 19
    24 = 11
    25 = 20
    26 = 15
    27 = 21
}

$cols = @("D", "J", "K", "L", "M", "P")

# First, snapshot the original values for every affected column/row so that
# writes do not clobber values that still need to be read as a source.
$orig = @{}
foreach ($col in $cols) {
    for ($r = 2; $r -le 27; $r++) {
        $orig["$col$r"] = $ws.Range("$col$r").Value2
    }
}

# Now write the permuted values back.
foreach ($col in $cols) {
    for ($r = 2; $r -le 27; $r++) {
        $srcRow = $rowMap[$r]
        $ws.Range("$col$r").Value = $orig["$col$srcRow"]
    }
}
